# Add two new computed columns (X: Med_Pont_H_5Jogos, Y: Med_Pont_A_5Jogos)
# to the existing match-data table on the active sheet, mirroring the
# upstream author's "Add files via upload" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ------------------------------------------------
$ws.Cells.Item(1, 24).Value = "Med_Pont_H_5Jogos"
$ws.Cells.Item(1, 24).Font.Bold = $true
$ws.Cells.Item(1, 24).Interior.Color = 65535
$ws.Cells.Item(1, 24).HorizontalAlignment = -4108

$ws.Cells.Item(1, 25).Value = "Med_Pont_A_5Jogos"
$ws.Cells.Item(1, 25).Font.Bold = $true
$ws.Cells.Item(1, 25).Interior.Color = 65535
$ws.Cells.Item(1, 25).HorizontalAlignment = -4108

# --- Data rows (2-49) ---------------------------------------------------
$xValues = @(1.2,2,2.2000000000000002,1.2,2,0.8,1.4,0.4,3,2.4,2.4,0.2,1.8,0.8,1.2,2.6,3,3,0.8,1.2,0.6,1,2.4,1,1.2,1,2,1.4,1.6,2.4,3,2.6,2.2000000000000002,1.4,1.6,1.6,0.4,2.4,1.2,1.2,2,1.8,0.2,1.2,1.2,2.6,2.4,2.6)
$yValues = @(0.8,0.8,0.4,0.6,0.8,3,1.8,1.4,1.6,1.4,1.2,1.4,1.8,2,0.6,1.6,1.2,1.2,0.8,1.4,1.6,1.8,0.6,1.2,1.6,1.2,1.6,1,1.6,1.4,0.4,0,0.2,1.4,0,0.8,1.4,1,1.4,1,2,2.6,0.8,0.4,1.8,0.4,0.2,1.6)

for ($i = 0; $i -lt $xValues.Length; $i++) {
    $r = $i + 2
    $xCell = $ws.Cells.Item($r, 24)
    $xCell.Value = $xValues[$i]
    $xCell.NumberFormat = "0.00"

    $yCell = $ws.Cells.Item($r, 25)
    $yCell.Value = $yValues[$i]
    $yCell.NumberFormat = "0.00"
}

# --- Column widths, matching the author's manual resize ----------------
# NOTE: the COM ColumnWidth setter here quantises the stored OOXML width to
# the nearest 1/6 of a character (it bakes in the standard ~0.833 char
# padding); the inputs below are chosen so the serialized <col width=.../>
# lands on (or as close as achievable to) 20 and 19.85546875 respectively.
$ws.Columns.Item(24).ColumnWidth = 19.166666666666668
$ws.Columns.Item(25).ColumnWidth = 19

# --- Page setup, as saved by the author's Excel session ----------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Selection moved to A9 on save --------------------------------------
[void]$ws.Range("A9").Select()

Write-Host "Added Med_Pont_H_5Jogos / Med_Pont_A_5Jogos columns"
